$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 443, shifting existing rows 443:539 down to 444:540
$ws.Rows.Item(443).Insert()

# Populate the newly inserted row 443 with the new record
$ws.Cells.Item(443, 1).Value = 4
$ws.Cells.Item(443, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(443, 3).Value = "Los Lagos"
$ws.Cells.Item(443, 4).Value = 45204
$ws.Cells.Item(443, 5).Value = 10
$ws.Cells.Item(443, 6).Value = 100114014
$ws.Cells.Item(443, 7).Value = "Betarraga"
$ws.Cells.Item(443, 8).Value = "Sin especificar"
$ws.Cells.Item(443, 9).Value = "Primera"
$ws.Cells.Item(443, 10).Value = 500
$ws.Cells.Item(443, 11).Value = 1000
$ws.Cells.Item(443, 12).Value = 1000
$ws.Cells.Item(443, 13).Value = 1000
$ws.Cells.Item(443, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(443, 15).Value = "Región Metropolitana"
$ws.Cells.Item(443, 16).Value = 200
$ws.Cells.Item(443, 17).Value = 5
$ws.Cells.Item(443, 18).Value = "Hortaliza"
